$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("privateinfo")
$ws.Range("A29").Value = "admin1"
